$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell: O1 = "notes" ---
$ws.Range("O1").Value = "notes"

# --- Row 2: fill in nchoices (L2) and package (N2) for the existing events/nothing event ---
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = "default"

# --- Row 3: flesh out the events/overheat event ---
# A3 (pathString) and C3 (baseprob) already hold the right values; only the
# description (E3) changes from "Overheat!" to the full event text.
$ws.Range("E3").Value = "The AI wakes to a wall of blinding brightness. The seedship's course has taken it close to a super-giant star that has proved to be far hotter than the guidance system anticipated, and the ship is gathering heat faster than the radiator fins can radiate it away. The AI must shut down part of the heat regulation system before the entire system catastrophically fails."
$ws.Range("E3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 152.55

# --- Row 4: new choice events/overheat/colonists ---
$ws.Range("A4").Value = "events/overheat/colonists"
$ws.Range("B4").Value = "Allow the sleep chambers to overheat"
$ws.Range("B4").WrapText = $true
$ws.Range("E4").Value = "The AI channels excess heat into the sleep chambers. The heat regulation system recovers as the ship moves away from the super-hot star, but not before [?-75] colonists have sustained tissue damage too severe for them to be successfully revived."
$ws.Range("E4").WrapText = $true
$ws.Range("L4").Value = 0
$ws.Rows.Item(4).RowHeight = 106.1

# --- Column widths: description column narrows now that it wraps ---
$ws.Columns.Item(5).ColumnWidth = 25.1

# --- Selection, matching the author's cursor position when the edit was saved ---
$ws.Range("I3").Select() | Out-Null
